# Edit: change the shared-string label in cell A1 of the "Demanda" sheet
# from "TX" to "T", and move the active selection to A2 (it was A5 before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demanda")

# Update the label text.
$ws.Range("A1").Value = "T"

# Restore/update the active selection on that sheet.
$ws.Activate()
$ws.Range("A2").Select()
